# Edit script: add 2022-Q4 data
# - Insert a new worksheet "2022-Q4" right after the "总计" (Total) summary sheet,
#   pushing all the existing quarterly sheets (2022-Q3 .. 2021-Q1) one position later
#   (their names/content are untouched).
# - Add a new summary row for 2022-Q4 at the top of the "总计" sheet's data table.

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q4" worksheet right after the summary sheet.
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $totalSheet)
$q4.Name = "2022-Q4"

# Header row (matches the other quarterly sheets)
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Data rows
$q4.Range("A2").Value = 0
$q4.Range("B2").NumberFormat = "@"
$q4.Range("B2").Value = "012977"
$q4.Range("C2").Value = "瑞达鑫红量化6个月持有混合A"
$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "0.35"
$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "94.66"
$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "4.94"
$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "0.0173"
$q4.Range("H2").Value = 5

$q4.Range("A3").Value = 1
$q4.Range("B3").NumberFormat = "@"
$q4.Range("B3").Value = "012978"
$q4.Range("C3").Value = "瑞达鑫红量化6个月持有混合C"
$q4.Range("D3").NumberFormat = "@"
$q4.Range("D3").Value = "0.09"
$q4.Range("E3").NumberFormat = "@"
$q4.Range("E3").Value = "94.66"
$q4.Range("F3").NumberFormat = "@"
$q4.Range("F3").Value = "4.94"
$q4.Range("G3").NumberFormat = "@"
$q4.Range("G3").Value = "0.0044"
$q4.Range("H3").Value = 5

# Apply the same formatting used on the other sheets' header row / index column
# (bold, thin border, centered) by copying it from the "总计" sheet's A2 cell,
# which already carries that exact style.
$totalSheet.Range("A2").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2:A3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Insert the new 2022-Q4 summary row into the "总计" sheet.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

# The freshly inserted row can pick up stray formatting (e.g. bold) - start
# from a clean slate for the plain (unstyled) B:D cells.
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.02

# Re-apply the index-column style to A2 (row insert does not carry it down).
# A3 still holds the original styled cell (it was A2 before the row insert).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

# Re-number the index column (A) for every data row so it stays a simple
# 0-based row counter, since the inserted row shifted the old values down.
for ($r = 2; $r -le 9; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

Write-Host "Done"
